$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New result for person #6 (row 7): 0.4179 -> 0.4333.
# B32 holds =AVERAGE(B2:B31), so it recalculates automatically.
$ws.Range("B7").Value = 0.4333

# Move the selection / view back to the top of the sheet (C1) instead of
# where it had been left (C32, with the window scrolled to row 6).
$ws.Range("C1").Select()
